$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "run with 10000 mutual information values" - the term list in column A got
# reshuffled and the mutual-information scores in column B were refreshed
# from the new run. Final state per row (A2:B11):
$labels = @("5", "9", "1", "3", "4", "6", "7", "8", "10", "2")
$values = @(5.0, 5.0, 4.0, 1.0, 2.0, 4.0, 4.0, 4.0, 1.0, 4.0)

for ($i = 0; $i -lt $labels.Length; $i++) {
    $row = $i + 2
    $existingLabel = [string]$ws.Cells.Item($row, 1).Value2
    if ($existingLabel -ne $labels[$i]) {
        # Leading apostrophe forces the numeric-looking term label to stay
        # text, matching the existing string-typed cells in column A.
        $ws.Cells.Item($row, 1).Formula = "'" + $labels[$i]
    }
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

$wb.Save()
